$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# Numeric-looking Price text values are prefixed with a literal leading apostrophe
# so Excel stores them as text (matching the sheet's inline-string convention)
# instead of auto-converting them to numbers.
$ws.Range('D2').Value = '61.533.12'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '3.378.37'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('D4').Value = '''0.997'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').Value = '''575.04'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').Value = '''138.01'
$ws.Range('E6').Value = '  +7.37%  '
$ws.Range('D7').Value = '''0.998'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').Value = '3.381.06'
$ws.Range('E8').Value = '  +2.01%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').Value = '''7.64'
$ws.Range('E10').Value = '  +4.60%  '
$ws.Range('D11').Value = '''0.125'
$ws.Range('E11').Value = '  +5.24%  '
$ws.Range('D12').Value = '''0.393'
$ws.Range('E12').Value = '  +4.54%  '
$ws.Range('D13').Value = '3.947.54'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('E15').Value = '  +5.05%  '
$ws.Range('D16').Value = '3.363.67'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').Value = '''25.48'
$ws.Range('E17').Value = '  +3.37%  '
$ws.Range('D18').Value = '61.352.68'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').Value = '''13.98'
$ws.Range('E19').Value = '  +4.93%  '
$ws.Range('E20').Value = '  +3.49%  '
$ws.Range('D21').Value = '''9.41'
$ws.Range('E21').Value = '  +4.75%  '
$ws.Range('D22').Value = '''381.81'
$ws.Range('E22').Value = '  +7.56%  '
$ws.Range('D23').Value = '''0.569'
$ws.Range('E23').Value = '  +2.04%  '
$ws.Range('D24').Value = '3.501.94'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = '''70.57'
$ws.Range('E26').Value = '  +1.60%  '
$ws.Range('D27').Value = '''0.0000123'
$ws.Range('E27').Value = '  +14.54%  '
$ws.Range('E28').Value = '  +11.92%  '
$ws.Range('D29').Value = '''7.84'
$ws.Range('E29').Value = '  +8.28%  '
$ws.Range('D30').Value = '''0.991'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').Value = '''8.27'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('D32').Value = '''0.158'
$ws.Range('E32').Value = '  +5.03%  '
$ws.Range('D33').Value = '''2.14'
$ws.Range('E33').Value = '  +1.34%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').Value = '3.401.54'
$ws.Range('E35').Value = '  +1.63%  '
$ws.Range('D36').Value = '''23.57'
$ws.Range('E36').Value = '  +4.50%  '
$ws.Range('D37').Value = '''5.53'
$ws.Range('E37').Value = '  +4.17%  '
$ws.Range('D38').Value = '''7.07'
$ws.Range('E38').Value = '  +3.51%  '
$ws.Range('D39').Value = '''1.55'
$ws.Range('E39').Value = '  +4.45%  '
$ws.Range('D40').Value = '''160.91'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').Value = '''0.0792'
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('E44').Value = '  +9.51%  '
$ws.Range('D45').Value = '''4.45'
$ws.Range('E45').Value = '  +1.69%  '
$ws.Range('D46').Value = '''0.769'
$ws.Range('E46').Value = '  +3.39%  '
$ws.Range('D47').Value = '''41.46'
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('D48').Value = '''23.29'
$ws.Range('E48').Value = '  +5.30%  '
$ws.Range('D49').Value = '''6.98'
$ws.Range('E49').Value = '  +3.89%  '
$ws.Range('D50').Value = '''22.91'
$ws.Range('E50').Value = '  +7.24%  '
$ws.Range('D51').Value = '2.340.67'
$ws.Range('E51').Value = '  +8.24%  '

# Rows 42 and 43 swap which coin they list (Stacks now ranks above FirstDigitalUSD),
# each with refreshed Price/Volume values.
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = "'1.74"
$ws.Range('E42').Value = '  +11.86%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'0.997"
$ws.Range('E43').Value = '  -0.50%  '
